# Update bank statement sample data (statement_1.xlsx) to the next
# generated period, per the commit "Update generated bank statements and
# ground truth CSV".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Account holder / card info ---
$ws.Range("C2").Value = "Hartmut"
# Card number is all digits; a leading apostrophe keeps it text (like the
# source data) instead of Excel auto-coercing it to a number.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 08.04.2025"

# --- Transaction rows 6-10 (dates / description / amount) ---
$ws.Range("B6").Value = "11.04."
$ws.Range("C6").Value = "12.04."
$ws.Range("D6").Value = "KARTENZ./11.04 EDEKA RO"
$ws.Range("E6").Value = "64,97-"

$ws.Range("B7").Value = "13.04."
$ws.Range("C7").Value = "14.04."
$ws.Range("D7").Value = "MCDONALDS Seelow"
$ws.Range("E7").Value = "35,81-"

$ws.Range("B8").Value = "16.04."
$ws.Range("C8").Value = "17.04."
$ws.Range("D8").Value = "BURGER KING Wismar"
$ws.Range("E8").Value = "16,75-"

$ws.Range("B9").Value = "18.04."
$ws.Range("C9").Value = "19.04."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "25,12-"

$ws.Range("B10").Value = "21.04."
$ws.Range("C10").Value = "22.04."
$ws.Range("D10").Value = "PAYPAL XUBVUF"
$ws.Range("E10").Value = "72,41-"

# --- Row 11 was a blank spacer row; it now holds a 6th transaction ---
$ws.Range("B11").Value = "23.04."
$ws.Range("C11").Value = "24.04."
$ws.Range("D11").Value = "BEITRAG Allianz SE K-29921652"
$ws.Range("E11").Value = "55,76-"
# Match the formatting used by the other amount cells in the column
# (right aligned, single line, not vertically centered) instead of the
# blank-row style it inherited before.
$ws.Range("E11").VerticalAlignment = -4107
$ws.Range("E11").WrapText = $false

# --- Closing balance / next statement date ---
$ws.Range("D12").Value = "KONTOSTAND AM 28.04.2025"
$ws.Range("E12").Value = "270,82-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 04.05.2025"
